$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.294.36"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.793.48"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.90"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.595"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.10"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.290"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0672"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0960"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "2.050.72"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.13"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").Value = "1.815.07"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.628"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "34.259.09"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.87"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("D35").Value = "1.356.43"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.640"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "80.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.930"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  +5.25%  "
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0497"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").Value = "1.952.81"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.75"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.04%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("E51").Value = "  -8.32%  "
